$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constraints")

$ws.Range("A3").Value  = "Created At date cannot be after current date"
$ws.Range("A6").Value  = "Created At date cannot be after current date"
$ws.Range("A11").Value = "Commencement date cannot be after current date"
$ws.Range("A12").Value = "Created At date cannot be after current date"
$ws.Range("A16").Value = "Created At date cannot be after current date"
$ws.Range("A17").Value = "Birth date cannot be after current date"
$ws.Range("A20").Value = "Created At date cannot be after current date"
$ws.Range("A21").Value = "Birth date cannot be after current date"
$ws.Range("A26").Value = "Created At date cannot be after current date"
$ws.Range("A29").Value = "Created At date cannot be after current date"
$ws.Range("A30").Value = "Amount must be a positive value"
$ws.Range("A31").Value = "Payment date cannot be after current date"
$ws.Range("A37").Value = "Created At date cannot be after current date"
$ws.Range("A40").Value = "Created At date cannot be after current date"
$ws.Range("A41").Value = "Email must have a valid format"
$ws.Range("A42").Value = "Type can only be C for Company or I for Individual"
$ws.Range("A44").Value = "Created At date cannot be after current date"

$ws.Columns.Item(1).ColumnWidth = 46.140625
